$wb = $excel.ActiveWorkbook

# The "Metadata" sheet holds Property/Value pairs in columns A/B.
$ws = $wb.Worksheets.Item("Metadata")

# Row 4 ("Name") currently has an empty Value cell; set it to "SexeVs".
$ws.Range("B4").Value = "SexeVs"

# Row 8 ("Date") Value cell needs its timestamp updated.
$ws.Range("B8").Value = "2025-07-18T06:40:38+00:00"
